$wb = $excel.ActiveWorkbook

# "0_First_176_Last_176" sheet: fill in the Status column (U) for rows 2,4,5,6,7
$ws1 = $wb.Worksheets.Item("0_First_176_Last_176")
$ws1.Range("U2").Value = "NA"
$ws1.Range("U4").Value = "Liked"
$ws1.Range("U5").Value = "Match"
$ws1.Range("U6").Value = "Email"
$ws1.Range("U7").Value = "Rejected"

# "1_First_336_Last_336" sheet: fill in the Status column (U) for row 2
$ws2 = $wb.Worksheets.Item("1_First_336_Last_336")
$ws2.Range("U2").Value = "NA"
